# Generate Report for Handoff
# Updates the localization-status report: mark several rows as high ("ht")
# priority, and refresh the "Latest HO Xliff Generate Date" / "Latest
# Handoff Datetime" timestamps for the corresponding rows.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 11, 12, 13, 14)

# Overview sheet: bump the "Latest HO Xliff Generate Date" column (G)
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-31 18:24:53"
}

# de-de sheet: mark Priority column (E) as "ht" and bump the
# "Latest Handoff Datetime" column (H) to match the Overview timestamp
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-31 18:24:53"
}

# zh-cn sheet: mark Priority column (E) as "ht" and bump the
# "Latest Handoff Datetime" column (H) to its own refreshed timestamp
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-31 18:24:48"
}
